$d = $word.ActiveDocument

# The document currently has 3 paragraphs:
#   1) "[add ideas]" (holds the _GoBack bookmark)
#   2) empty (indented)
#   3) empty (indented)
# Replace the whole body's paragraph range (1..3) with the full set of
# project-idea paragraphs from the commit, re-attaching the _GoBack
# bookmark to a final bare paragraph.

$lastPara = $d.Paragraphs($d.Paragraphs.Count)
$target = $d.Range(0, $lastPara.Range.End)

$newBodyXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:ind w:left="720"/></w:pPr><w:r><w:t>Smart Plant Watering System: Use the soil moisture sensor to monitor the moisture level in a plant pot and activate the water pump when the soil becomes too dry. You can also incorporate the temperature sensor to adjust the watering schedule based on ambient conditions.</w:t></w:r></w:p><w:p><w:pPr><w:ind w:left="720"/></w:pPr></w:p><w:p><w:pPr><w:ind w:left="720"/></w:pPr><w:r><w:t>Intruder Alert System: Utilize the motion sensor and the buzzer to create an alarm system. When motion is detected, the buzzer will sound an alarm, and you can add LED lights for visual indication. You can even integrate the camera module to capture images or videos of the intruder.</w:t></w:r></w:p><w:p><w:pPr><w:ind w:left="720"/></w:pPr></w:p><w:p><w:pPr><w:ind w:left="720"/></w:pPr><w:r><w:t>Home Weather Station: Connect the temperature sensor, light sensor, and distance sensor to collect weather data. Display the temperature, light intensity, and proximity information on an LCD screen or a web interface. You can even store the data in a database for future analysis.</w:t></w:r></w:p><w:p><w:pPr><w:ind w:left="720"/></w:pPr></w:p><w:p><w:pPr><w:ind w:left="720"/></w:pPr><w:r><w:t>Smart Home Lighting Control: Use the light sensor to automatically control LED lights or lamps based on ambient light levels. For example, turn on the lights when it gets dark and turn them off when it becomes bright. You can also incorporate a joystick or a slider to adjust the brightness.</w:t></w:r></w:p><w:p/><w:p><w:pPr><w:ind w:left="720"/></w:pPr><w:r><w:t>Aquarium Monitoring System: Use the water depth sensor and the temperature sensor to monitor the water level and temperature in an aquarium. Display the data on an LCD screen or send notifications when certain conditions are met, such as low water level or high temperature.</w:t></w:r></w:p><w:p/><w:p><w:pPr><w:ind w:left="720"/></w:pPr><w:r><w:t>Automated Pet Feeder: Create an automated pet feeding system using the servo motor, push buttons, and LCD screen. Set up a schedule for feeding your pet, and when the designated time arrives, the servo motor dispenses the food. The LCD screen can display the feeding schedule and other relevant information.</w:t></w:r></w:p><w:p/><w:p><w:pPr><w:ind w:left="720"/></w:pPr><w:r><w:t>Garage Door Opener: Connect a relay to the Raspberry Pi and wire it to your garage door opener system. Use a push button or a motion sensor to trigger the relay and open/close the garage door. You can even add a security feature by integrating a camera module to monitor the garage entrance.</w:t></w:r></w:p><w:p><w:pPr><w:ind w:left="720"/></w:pPr></w:p><w:p><w:pPr><w:ind w:left="720"/></w:pPr><w:r><w:t>Weather-Activated Window Control: Utilize the temperature sensor and distance sensor to automate your windows based on weather conditions. For example, when the temperature rises above a certain threshold, the windows can automatically open. The distance sensor can prevent the windows from opening if an obstacle is detected.</w:t></w:r></w:p><w:p/><w:p><w:pPr><w:ind w:left="720"/></w:pPr><w:r><w:t xml:space="preserve">Personal Weather Station: Set up a personal weather station using the temperature sensor, humidity </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>sensor</w:t></w:r><w:r><w:t>(</w:t></w:r><w:proofErr w:type="spellStart"/><w:proofErr w:type="gramEnd"/><w:r><w:t>Mr</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> Brake has one you could use)</w:t></w:r><w:r><w:t>. Collect weather data, such as temperature, humidity, and display it on an LCD screen or a web interface. You can even log the data and create weather graphs over time.</w:t></w:r></w:p><w:p><w:pPr><w:ind w:left="720"/></w:pPr></w:p><w:p><w:pPr><w:ind w:left="720"/></w:pPr><w:r><w:t>Automated Curtain Control: Connect a stepper motor or a servo motor to your curtain system. Use a light sensor or a motion sensor to automatically open or close the curtains based on ambient light levels or the presence of people in the room.</w:t></w:r></w:p><w:p><w:pPr><w:ind w:left="720"/></w:pPr></w:p><w:p><w:pPr><w:ind w:left="720"/></w:pPr><w:r><w:t>Automated Plant Monitoring: Utilize the light sensor, temperature sensor, and soil moisture sensor to monitor the conditions for your plants. Collect data on light intensity, temperature, and soil moisture levels. Set up notifications or display the data on an LCD screen to remind you when to water or adjust lighting for optimal plant growth.</w:t></w:r></w:p><w:p><w:pPr><w:ind w:left="720"/></w:pPr></w:p><w:p><w:pPr><w:ind w:left="720"/></w:pPr><w:r><w:t>Smart Home Thermostat: Use the temperature sensor and LCD screen to create a smart thermostat for your home. Program the Raspberry Pi to monitor the temperature and control your heating or cooling system based on predefined temperature thresholds. You can also incorporate a joystick or a slider to adjust the desired temperature manually.</w:t></w:r></w:p><w:p><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p>'

$target.InsertXML($newBodyXml)

Write-Output "Paragraphs after edit: $($d.Paragraphs.Count)"
